# BIIBBag.xlsx update:
#  - Row 3 gets a PriceChange (X3) and an UpDown verdict of "Up" (Y3),
#    since "moved against" was creating false positives and is no longer
#    left blank for this day.
#  - A brand new observation row (row 4) is appended with the day's
#    sentiment / fundamentals / technical data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- complete row 3 (previously missing PriceChange/UpDown) ---
$ws.Range("X3").Value = 4.75
$ws.Range("Y3").Value = "Up"

# --- new row 4 ---
# Copy A3's format (date/time number format) down to A4 before setting the value
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 42641.891493055555

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = -4
$ws.Range("E4").Value = 14968
$ws.Range("F4").Value = 1683
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 42
$ws.Range("J4").Value = 57
$ws.Range("K4").Value = 8103
$ws.Range("L4").Value = 207
$ws.Range("M4").Value = 129
$ws.Range("N4").Value = 12
$ws.Range("O4").Value = 16
$ws.Range("P4").Value = "Bag"
$ws.Range("Q4").Value = 57.519894101767122
$ws.Range("R4").Value = 1.83

# Copy S3:T3's format (percentage number format) down to S4:T4 before setting values
$ws.Range("S3:T3").Copy($ws.Range("S4:T4"))
$ws.Range("S4").Value = 0.13639999999999999
$ws.Range("T4").Value = 0.0165

$ws.Range("U4").Value = 6.04
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 2
